$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-11) from 2 to 4
$ws.Range("A2:A11").Value = 4

# Update the selection/active cell on the sheet view
$ws.Range("A2:A11").Select()
